# Auto-generated edit script: updates cryptos worksheet price/volume data
# per commit 'Updated cryptos list on Wed Jun 19 05:47:11 UTC 2024 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.466.43'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '3.557.50'
$ws.Range("E3").Value = '  +2.92%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''601.44'
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("D6").Value = '''140.53'
$ws.Range("E6").Value = '  +2.35%  '
$ws.Range("D7").Value = '3.558.76'
$ws.Range("E7").Value = '  +3.04%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = '''0.493'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").Value = '''0.125'
$ws.Range("E10").Value = '  +2.79%  '
$ws.Range("D11").Value = '''7.04'
$ws.Range("E11").Value = '  -6.06%  '
$ws.Range("D12").Value = '''0.395'
$ws.Range("E12").Value = '  +3.66%  '
$ws.Range("D13").Value = '4.166.90'
$ws.Range("E13").Value = '  +3.12%  '
$ws.Range("D14").Value = '''0.0000187'
$ws.Range("E14").Value = '  +2.30%  '
$ws.Range("D15").Value = '''27.16'
$ws.Range("E15").Value = '  +1.61%  '
$ws.Range("D16").Value = '3.567.37'
$ws.Range("E16").Value = '  +3.32%  '
$ws.Range("E17").Value = '  +1.36%  '
$ws.Range("D18").Value = '65.530.38'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").Value = '''10.28'
$ws.Range("E19").Value = '  +3.35%  '
$ws.Range("D20").Value = '''5.90'
$ws.Range("E20").Value = '  +1.95%  '
$ws.Range("D21").Value = '''14.27'
$ws.Range("E21").Value = '  +3.40%  '
$ws.Range("D22").Value = '''396.25'
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = '''0.573'
$ws.Range("E23").Value = '  +4.14%  '
$ws.Range("D24").Value = '3.707.11'
$ws.Range("E24").Value = '  +2.92%  '
$ws.Range("D25").Value = '''74.35'
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("D27").Value = '''0.0000116'
$ws.Range("E27").Value = '  +8.95%  '
$ws.Range("D28").Value = '''7.91'
$ws.Range("E28").Value = '  +9.36%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("D31").Value = '''8.34'
$ws.Range("E31").Value = '  +0.70%  '
$ws.Range("D32").Value = '3.579.90'
$ws.Range("E32").Value = '  +3.44%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '''1.00'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  +0.76%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '''23.91'
$ws.Range("E35").Value = '  +4.11%  '
$ws.Range("D36").Value = '''1.27'
$ws.Range("E36").Value = '  +3.83%  '
$ws.Range("D37").Value = '''7.09'
$ws.Range("E37").Value = '  +2.29%  '
$ws.Range("E38").Value = '  +1.09%  '
$ws.Range("D39").Value = '''167.98'
$ws.Range("E39").Value = '  -2.93%  '
$ws.Range("E40").Value = '  +4.51%  '
$ws.Range("D41").Value = '''0.0806'
$ws.Range("E41").Value = '  +3.27%  '
$ws.Range("D42").Value = '''0.834'
$ws.Range("E42").Value = '  +1.31%  '
$ws.Range("D43").Value = '''26.74'
$ws.Range("E43").Value = '  +14.69%  '
$ws.Range("D44").Value = '''43.02'
$ws.Range("E44").Value = '  -1.32%  '
$ws.Range("D45").Value = '''1.00'
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = '''4.45'
$ws.Range("E46").Value = '  -0.34%  '
$ws.Range("E47").Value = '  +3.58%  '
$ws.Range("E48").Value = '  +8.21%  '
$ws.Range("D49").Value = '2.451.72'
$ws.Range("E49").Value = '  +10.33%  '
$ws.Range("D50").Value = '''6.83'
$ws.Range("E50").Value = '  +3.71%  '
$ws.Range("B51").Value = 'LidoDAOToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D51").Value = '''2.38'
$ws.Range("E51").Value = '  +20.47%  '
